$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.720326527486293
$ws.Range("C2").Value = 0.4125059731021565
$ws.Range("E2").Value = 0.05134008465251938
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.002460635343412849
$ws.Range("M2").Value = 0.5172146108966018
$ws.Range("N2").Value = 1.808380866996032
$ws.Range("B3").Value = 1.560570256317419
$ws.Range("C3").Value = 0.364136531653628
$ws.Range("E3").Value = 0.04819871944929943
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.002467433184192771
$ws.Range("M3").Value = 0.4695452906511406
$ws.Range("N3").Value = 1.798734719878581
$ws.Range("B4").Value = 1.46350614791703
$ws.Range("C4").Value = 0.3346175515906111
$ws.Range("E4").Value = 0.04630699582197551
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.002471816524338695
$ws.Range("M4").Value = 0.4405957011175516
$ws.Range("N4").Value = 1.793352692763079
$ws.Range("B5").Value = 1.424203608197104
$ws.Range("C5").Value = 0.3226314632369736
$ws.Range("E5").Value = 0.04554525756180539
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002473655651067034
$ws.Range("M5").Value = 0.4288767996643799
$ws.Range("N5").Value = 1.791293632154677
$ws.Range("B6").Value = 1.417692480519804
$ws.Range("C6").Value = 0.3206437257967139
$ws.Range("E6").Value = 0.04541931945539091
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002473964236792288
$ws.Range("M6").Value = 0.4269355538883346
$ws.Range("N6").Value = 1.790959776780284
$ws.Range("B7").Value = 1.462975088607891
$ws.Range("C7").Value = 0.3344557311616256
$ws.Range("E7").Value = 0.04629668591170955
$ws.Range("F7").Value = 0.3529483938368969
$ws.Range("G7").Value = 0.002471841113047194
$ws.Range("M7").Value = 0.4404373416213971
$ws.Range("N7").Value = 1.793324382730404
$ws.Range("B8").Value = 1.665025534308086
$ws.Range("C8").Value = 0.3957895082122604
$ws.Range("E8").Value = 0.05024913946652632
$ws.Range("F8").Value = 0.4248636149813905
$ws.Range("G8").Value = 0.002462935899910131
$ws.Range("M8").Value = 0.5007106354897672
$ws.Range("N8").Value = 1.804941541416397
$ws.Range("B9").Value = 2.06970441586418
$ws.Range("C9").Value = 0.517596977161304
$ws.Range("E9").Value = 0.05830241937369252
$ws.Range("F9").Value = 0.5661985755042025
$ws.Range("G9").Value = 0.002447124736167261
$ws.Range("M9").Value = 0.6215424825642799
$ws.Range("N9").Value = 1.832097583221767
$ws.Range("B10").Value = 2.372644246010225
$ws.Range("C10").Value = 0.6081763660110369
$ws.Range("E10").Value = 0.06441580084739229
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002436501448746409
$ws.Range("M10").Value = 0.7120744492649607
$ws.Range("N10").Value = 1.854836238832547
$ws.Range("B11").Value = 2.51178494891343
$ws.Range("C11").Value = 0.6496537634676542
$ws.Range("E11").Value = 0.06724236942517337
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002431881309022619
$ws.Range("M11").Value = 0.7536748774556941
$ws.Range("N11").Value = 1.865812329107342
$ws.Range("B12").Value = 2.564673536994576
$ws.Range("C12").Value = 0.6654020271670333
$ws.Range("E12").Value = 0.06831948445970681
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002430162100254495
$ws.Range("M12").Value = 0.7694904564404794
$ws.Range("N12").Value = 1.87006174226326
$ws.Range("B13").Value = 2.553274060268848
$ws.Range("C13").Value = 0.6620084628902987
$ws.Range("E13").Value = 0.06808720475554253
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002430531016684722
$ws.Range("M13").Value = 0.7660814738765538
$ws.Range("N13").Value = 1.869142384803723
$ws.Range("B14").Value = 2.516132092546059
$ws.Range("C14").Value = 0.6509485328502365
$ws.Range("E14").Value = 0.06733084776887921
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.00243173926197319
$ws.Range("M14").Value = 0.7549747699102056
$ws.Range("N14").Value = 1.866160054752328
$ws.Range("B15").Value = 2.493407709323037
$ws.Range("C15").Value = 0.6441795087798141
$ws.Range("E15").Value = 0.06686844306963025
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002432483291544987
$ws.Range("M15").Value = 0.74817979262194
$ws.Range("N15").Value = 1.864345466545871
$ws.Range("B16").Value = 2.363578511735
$ws.Range("C16").Value = 0.6054714144522677
$ws.Range("E16").Value = 0.06423201201464224
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002436807642710526
$ws.Range("M16").Value = 0.7093643613625886
$ws.Range("N16").Value = 1.854131825153161
$ws.Range("B17").Value = 2.28427890926605
$ws.Range("C17").Value = 0.581796780794491
$ws.Range("E17").Value = 0.06262645934910438
$ws.Range("F17").Value = 0.6400460337215605
$ws.Range("G17").Value = 0.002439514754064557
$ws.Range("M17").Value = 0.6856608681564751
$ws.Range("N17").Value = 1.848029483918197
$ws.Range("B18").Value = 2.238792776550667
$ws.Range("C18").Value = 0.5682052421475419
$ws.Range("E18").Value = 0.0617072638762366
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002441091822334767
$ws.Range("M18").Value = 0.6720663450315669
$ws.Range("N18").Value = 1.844578953066332
$ws.Range("B19").Value = 2.223413179643501
$ws.Range("C19").Value = 0.5636076822340783
$ws.Range("E19").Value = 0.06139676836360408
$ws.Range("F19").Value = 0.6191636801734006
$ws.Range("G19").Value = 0.002441629233571319
$ws.Range("M19").Value = 0.6674701150336375
$ws.Range("N19").Value = 1.843420796690538
$ws.Range("B20").Value = 2.292707506711054
$ws.Range("C20").Value = 0.5843143304404066
$ws.Range("E20").Value = 0.06279692916800528
$ws.Range("F20").Value = 0.6429339538360921
$ws.Range("G20").Value = 0.00243922450822489
$ws.Range("M20").Value = 0.6881800813649761
$ws.Range("N20").Value = 1.848672928415084
$ws.Range("B21").Value = 2.527036126684436
$ws.Range("C21").Value = 0.6541959511528717
$ws.Range("E21").Value = 0.06755282318246714
$ws.Range("F21").Value = 0.7228739723492197
$ws.Range("G21").Value = 0.002431383549404309
$ws.Range("M21").Value = 0.7582353648812443
$ws.Range("N21").Value = 1.86703349630659
$ws.Range("B22").Value = 2.681347096063803
$ws.Range("C22").Value = 0.7001117373686725
$ws.Range("E22").Value = 0.07070054498927192
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002426435762958231
$ws.Range("M22").Value = 0.804385426887535
$ws.Range("N22").Value = 1.879576298125272
$ws.Range("B23").Value = 2.598879426504027
$ws.Range("C23").Value = 0.6755824399731409
$ws.Range("E23").Value = 0.06901686672943441
$ws.Range("F23").Value = 0.7472568307916134
$ws.Range("G23").Value = 0.002429060388560833
$ws.Range("M23").Value = 0.7797200632110446
$ws.Range("N23").Value = 1.872831571003502
$ws.Range("B24").Value = 2.288896613608642
$ws.Range("C24").Value = 0.5831760862460555
$ws.Range("E24").Value = 0.06271984776090633
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002439355663834891
$ws.Range("M24").Value = 0.6870410426847826
$ws.Range("N24").Value = 1.848381847324987
$ws.Range("B25").Value = 1.959274939544571
$ws.Range("C25").Value = 0.4844651254614405
$ws.Range("E25").Value = 0.05609015891679547
$ws.Range("F25").Value = 0.5279251897347308
$ws.Range("G25").Value = 0.002451226659085448
$ws.Range("M25").Value = 0.5885567523550463
$ws.Range("N25").Value = 1.82427049017717
